$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6 (previously "QA Engineer" row) to the new "Service Now Developer" posting.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Service Now Developer"
$ws.Range("C6").Value = "Chennai"

$description = "Test and deliver quality products.Mange Service now taksRequirements:`n- Bachelor’s or Master’s degree in Computer Science, Engineering, or a related field.`n- 3+ years of experience developing backend applications using Java and Spring Boot.`n- Strong understanding of RESTful API design, HTTP protocol, and stateless architecture.`n- Experience with databases like MongoDB, MySQL, or PostgreSQL.`n- Familiarity with containerization technologies such as Docker and orchestration tools like Kubernetes is a plus.`n- Understanding of software design patterns, data structures, and algorithms.`n- Strong problem-solving and analytical thinking skills.`n- Excellent written and verbal communication skills."

$ws.Range("D6").Value = $description
$ws.Range("D6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 135
